$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 18: "Fan Duct" section header (bold, like the other section headers) ---
$ws.Range("A18").Value = "Fan Duct"
$ws.Range("A18").Font.Bold = $true

# --- Row 19: M3x30 Male Standoff ---
$ws.Range("A19").Value = "M3x30 Male Standoff"
$ws.Range("B19").Value = 4

# --- Row 20: M3x40 Female Standoff ---
$ws.Range("A20").Value = "M3x40 Female Standoff"
$ws.Range("B20").Value = 4

# --- Row 21: M4x14 Cap Head Screw ---
$ws.Range("A21").Value = "M4x14 Cap Head Screw"
$ws.Range("B21").Value = 8
$ws.Range("D21").Value = "Depends on thickness of dry box walls"

# --- Row 22: 6020 Fan ---
$ws.Range("A22").Value = "6020 Fan"
$ws.Range("B22").Value = 1
$ws.Range("D22").Value = "High CFM"

# --- Row 23: PTC Heater 12x7cm 260C (+ hyperlink note) ---
$ws.Range("A23").Value = "PTC Heater 12x7cm 260C"
$ws.Range("B23").Value = 1
$ws.Hyperlinks.Add($ws.Range("E23"), "https://www.aliexpress.com/item/1005003758412972.html") | Out-Null
$ws.Range("D23").Value = "YEKMLCO Soldering Plate with Cord "

# --- Row 24: SSR ---
$ws.Range("A24").Value = "SSR"
$ws.Range("B24").Value = 1
$ws.Range("D24").Value = "10A for 120v service, 5a for 220v service"

# --- Row 25: Thermal Fuse 105C ---
$ws.Range("A25").Value = "Thermal Fuse 105C"
$ws.Range("B25").Value = 1

# Column D is now wider because of the new, longer notes in it
# (38.5 is the input that this host's column-width rounding resolves
# closest to Excel's real best-fit result of 39.28515625 characters)
$ws.Columns("D").ColumnWidth = 38.5

# Move the active selection past the newly entered data, like Excel does
# after typing values down a column and landing on the next empty row
$ws.Range("A26").Select() | Out-Null
